$wb = $excel.ActiveWorkbook

$wsBilling = $wb.Worksheets.Item("Billing")
$wsShipping = $wb.Worksheets.Item("Shipping")

# --- Shipping sheet: add an "id" column at the front (mirrors Billing's id
# column) and correct the postcode value.
$wsShipping.Columns.Item(1).Insert() | Out-Null
$wsShipping.Range("L2").Clear() | Out-Null
$wsShipping.Range("A1").Value = "id"
$wsShipping.Range("A2").Value = 155
$wsShipping.Range("I2").Value = 33121

# --- View-state: Shipping becomes the selected/active tab; Billing loses
# tabSelected and its selection moves.
$wsBilling.Activate() | Out-Null
$wsBilling.Range("A24").Select() | Out-Null

$wsShipping.Activate() | Out-Null
$wsShipping.Range("G8").Select() | Out-Null
